$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)

# --- Step 1: split the first paragraph's run into two runs ---
# "The Dog House should invest in advertising on the following platforms:"
# becomes "The Dog House should invest in advertising on the following " + "platforms:"
$tr = $shp.TextFrame.TextRange
$tail = $tr.Characters(61, 10)
$tail.Text = $tail.Text

# --- Step 2: append two new paragraphs after "Instagram" ---
$cr = [char]13
$tr2 = $shp.TextFrame.TextRange
$insertText = "" + $cr + $cr + "Additionally, it is recommended to add an Instagram button to the website along with the Facebook button."
$new = $tr2.InsertAfter($insertText)

# --- Step 3: suppress the bullet on the two new paragraphs (match "buNone" styling used elsewhere in this placeholder) ---
$tr3 = $shp.TextFrame.TextRange
$blankPara = $tr3.Characters(92, 0)
$blankPara.ParagraphFormat.Bullet.Visible = 0
$newTextPara = $tr3.Characters(93, 105)
$newTextPara.ParagraphFormat.Bullet.Visible = 0
